$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Append a new row to the (only) table in the document.
# ------------------------------------------------------------------
$table = $d.Tables.Item(1)
$newRow = $table.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "10"
$newRow.Cells.Item(2).Range.Text = "Critical"
$newRow.Cells.Item(3).Range.Text = "4h"
$newRow.Cells.Item(4).Range.Text = "As a supervisor I want documentation of the client/server system"

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark so that it sits in the last (empty)
#    paragraph of the document, right after the table, instead of its
#    old spot in the "23.04.2018" line - mirroring what Word itself
#    does to mark the most-recent edit location once the new row
#    above has been typed in.
#
#    The bookmark-placement call in this host only anchors reliably
#    against a range that currently contains text, so a temporary
#    placeholder character is typed into the empty trailing
#    paragraph, the bookmark is anchored around that placeholder, and
#    then the placeholder is removed again - leaving the (now empty)
#    bookmark sitting exactly where it should, and replacing the
#    previous "_GoBack" bookmark (bookmark names are unique).
# ------------------------------------------------------------------
$end = $d.Content.End
$sel = $word.Selection
$sel.SetRange($end - 1, $end)
$sel.TypeText("X")

$placeholderStart = $end - 1
$placeholder = $d.Range($placeholderStart, $placeholderStart + 1)
$d.Bookmarks.Add("_GoBack", $placeholder)

$placeholder = $d.Range($placeholderStart, $placeholderStart + 1)
$placeholder.Text = ""
